# Update the daily USD Amount figure in T2 and move the active selection
# to T3, matching the day's data refresh (upload of 15072025 data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the USD Amount value for the Roobic crypto deposit row.
$ws.Range("T2").Value = 294167

# Leave the selection where the user last clicked after the update.
$ws.Range("T3").Select()
